$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2594945728778839
$ws.Range("B1").Value = 0.3570446968078613
$ws.Range("C1").Value = 3.54664134979248
$ws.Range("D1").Value = 3.954626560211182
$ws.Range("E1").Value = 1.278937816619873
